$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume-change (column E) values.
# A leading apostrophe forces Excel to store the value as literal text (matching
# the source data, which keeps decimal-look-alike numbers such as "57.127.00" as text)
# and Style is reset to "Normal" afterwards so no stray number-format/style is left on the cell.

$ws.Range("D2").Value = "'57.127.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.59%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.316.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'533.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'132.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.69%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.24%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.339.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.47%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.43%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.18%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.344"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.25%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'23.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.14%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.736.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'57.178.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.45%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.333.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'339.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.21%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +2.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.17%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'61.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.48%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'170.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.21%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.55%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0722"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.82%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.51%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'18.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.89%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -3.24%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.904"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'148.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.69%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.90%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'278.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.56%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0929"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0504"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.557"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.62%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'18.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0216"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.52%  "
$ws.Range("E51").Style = "Normal"
